# Update the Year-End Reconciliation report template from the
# 2022-2023 school year to the 2024-2025 school year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value  = "RECONCILIATION REPORT FOR THE 2024-2025 SCHOOL YEAR"
$ws.Range("H9").Value  = "Total Amount Due for 2024-2025 School Year"
$ws.Range("G12").Value = "          Total Amount Due for 2024-2025 School Year:"
$ws.Range("C16").Value = "July, 2024"
$ws.Range("C22").Value = "January, 2025"
$ws.Range("G30").Value = "            Total Paid to Date for 2024-2025 School Year:"
